$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header rich-text runs: bulletin number and date range ---
$volCell = $ws.Range("A8")
$volCell.Characters(21, 2).Text = "49"

$dateCell = $ws.Range("C9")
$dateCell.Characters(27, 10).Text = "12/2/2024"
$dateCell.Characters(47, 9).Text = "12/8/2024"

# --- Weekly crime-statistics table updates ---
$ws.Range("D14").Value = "'0"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "***.*"
$ws.Range("E14").NumberFormat = "General"
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 75
$ws.Range("L15").Value = 27.272727272727
$ws.Range("M15").Value = 133.333333333333
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = "'0"
$ws.Range("C16").NumberFormat = "General"
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = -62.5
$ws.Range("J16").Value = 68
$ws.Range("K16").Value = -2.941176470588
$ws.Range("L16").Value = 29.411764705882
$ws.Range("M16").Value = -32.653061224489
$ws.Range("N16").Value = -87.54716981132
$ws.Range("C17").Value = "'0"
$ws.Range("C17").NumberFormat = "General"
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 71.428571428571
$ws.Range("J17").Value = 128
$ws.Range("K17").Value = -7.8125
$ws.Range("L17").Value = 18
$ws.Range("M17").Value = 18
$ws.Range("N17").Value = -57.090909090909
$ws.Range("C18").Value = "'0"
$ws.Range("C18").NumberFormat = "General"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = -6.818181818181
$ws.Range("L18").Value = -12.765957446808
$ws.Range("M18").Value = -64.655172413793
$ws.Range("N18").Value = -92.400370713623
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 20
$ws.Range("F19").Value = 22
$ws.Range("H19").Value = -24.137931034482
$ws.Range("I19").Value = 349
$ws.Range("J19").Value = 407
$ws.Range("K19").Value = -14.250614250614
$ws.Range("L19").Value = -28.775510204081
$ws.Range("M19").Value = 14.802631578947
$ws.Range("N19").Value = -21.0407239819
$ws.Range("C20").Value = 3
$ws.Range("E20").Value = 200
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 177
$ws.Range("J20").Value = 119
$ws.Range("K20").Value = 48.739495798319
$ws.Range("L20").Value = 80.612244897959
$ws.Range("M20").Value = 25.531914893617
$ws.Range("N20").Value = -90.390879478827
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 12
$ws.Range("E21").Value = -16.666666666666
$ws.Range("F21").Value = 48
$ws.Range("G21").Value = 53
$ws.Range("H21").Value = -9.43396226415
$ws.Range("I21").Value = 807
$ws.Range("J21").Value = 821
$ws.Range("K21").Value = -1.705237515225
$ws.Range("L21").Value = -4.609929078014
$ws.Range("M21").Value = -8.503401360544
$ws.Range("N21").Value = -80.785714285714
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 116
$ws.Range("G24").Value = 118
$ws.Range("H24").Value = -1.694915254237
$ws.Range("I24").Value = 1344
$ws.Range("J24").Value = 1456
$ws.Range("K24").Value = -7.692307692307
$ws.Range("L24").Value = -24.409448818897
$ws.Range("M24").Value = 18.518518518518
$ws.Range("C25").Value = 23
$ws.Range("D25").Value = 31
$ws.Range("E25").Value = -25.806451612903
$ws.Range("F25").Value = 62
$ws.Range("G25").Value = 80
$ws.Range("H25").Value = -22.5
$ws.Range("I25").Value = 816
$ws.Range("J25").Value = 799
$ws.Range("K25").Value = 2.127659574468
$ws.Range("L25").Value = -24.723247232472
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 30
$ws.Range("F26").Value = 39
$ws.Range("G26").Value = 39
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 351
$ws.Range("J26").Value = 363
$ws.Range("K26").Value = -3.305785123966
$ws.Range("L26").Value = 18.581081081081
$ws.Range("M26").Value = -6.14973262032
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 18
$ws.Range("K27").Value = 28.571428571428
$ws.Range("L27").Value = 20
$ws.Range("C28").Value = "'0"
$ws.Range("C28").NumberFormat = "General"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 50
$ws.Range("J28").Value = 28
$ws.Range("K28").Value = 71.428571428571
$ws.Range("F31").Value = 1
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("G31").Value = "'0"
$ws.Range("G31").NumberFormat = "General"
$ws.Range("H31").Value = "***.*"
$ws.Range("H31").NumberFormat = "General"
